$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated natmi (ligand-receptor) table now has 4 sending clusters
# (ECs, FAPs, M2, sCs) instead of 3, giving 16 data rows instead of 12.
# Clear the old A2:T13 block first so it can be rewritten at its new size (A2:T17).
$ws.Range("A2:T13").ClearContents()

# Row 2: ECs -> ECs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "F2"
$ws.Cells.Item(2, 3).Value = "Gp1ba"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 0.6740386666666667
$ws.Cells.Item(2, 8).Value = 2.022116
$ws.Cells.Item(2, 9).Value = 0.3058319194585966
$ws.Cells.Item(2, 10).Value = 0.3058319194585966
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 2.642196
$ws.Cells.Item(2, 14).Value = 7.926588000000001
$ws.Cells.Item(2, 15).Value = 0.26568831615543
$ws.Cells.Item(2, 16).Value = 0.26568831615543
$ws.Cells.Item(2, 17).Value = 1.780942268912
$ws.Cells.Item(2, 18).Value = 16.028480420208
$ws.Cells.Item(2, 19).Value = 0.08125596770753761
$ws.Cells.Item(2, 20).Value = 0.08125596770753764

# Row 3: ECs -> FAPs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "F2"
$ws.Cells.Item(3, 3).Value = "Gp1ba"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 0.6740386666666667
$ws.Cells.Item(3, 8).Value = 2.022116
$ws.Cells.Item(3, 9).Value = 0.3058319194585966
$ws.Cells.Item(3, 10).Value = 0.3058319194585966
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 4.666004666666667
$ws.Cells.Item(3, 14).Value = 13.998014
$ws.Cells.Item(3, 15).Value = 0.469194156323015
$ws.Cells.Item(3, 16).Value = 0.4691941563230151
$ws.Cells.Item(3, 17).Value = 3.145067564180445
$ws.Cells.Item(3, 18).Value = 28.305608077624
$ws.Cells.Item(3, 19).Value = 0.1434945494270245
$ws.Cells.Item(3, 20).Value = 0.1434945494270245

# Row 4: ECs -> M2 (ligand F2, receptor Gp1ba)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "F2"
$ws.Cells.Item(4, 3).Value = "Gp1ba"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 0.6740386666666667
$ws.Cells.Item(4, 8).Value = 2.022116
$ws.Cells.Item(4, 9).Value = 0.3058319194585966
$ws.Cells.Item(4, 10).Value = 0.3058319194585966
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 1.114591666666667
$ws.Cells.Item(4, 14).Value = 3.343775
$ws.Cells.Item(4, 15).Value = 0.1120787341732184
$ws.Cells.Item(4, 16).Value = 0.1120787341732184
$ws.Cells.Item(4, 17).Value = 0.7512778808777778
$ws.Cells.Item(4, 18).Value = 6.7615009279
$ws.Cells.Item(4, 19).Value = 0.03427725440268519
$ws.Cells.Item(4, 20).Value = 0.03427725440268519

# Row 5: ECs -> sCs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "F2"
$ws.Cells.Item(5, 3).Value = "Gp1ba"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 0.6740386666666667
$ws.Cells.Item(5, 8).Value = 2.022116
$ws.Cells.Item(5, 9).Value = 0.3058319194585966
$ws.Cells.Item(5, 10).Value = 0.3058319194585966
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 1.521928
$ws.Cells.Item(5, 14).Value = 4.565784000000001
$ws.Cells.Item(5, 15).Value = 0.1530387933483365
$ws.Cells.Item(5, 16).Value = 0.1530387933483365
$ws.Cells.Item(5, 17).Value = 1.025838319882667
$ws.Cells.Item(5, 18).Value = 9.232544878944001
$ws.Cells.Item(5, 19).Value = 0.04680414792134926
$ws.Cells.Item(5, 20).Value = 0.04680414792134926

# Row 6: FAPs -> ECs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "F2"
$ws.Cells.Item(6, 3).Value = "Gp1ba"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 0.9731926666666667
$ws.Cells.Item(6, 8).Value = 2.919578
$ws.Cells.Item(6, 9).Value = 0.4415672215387696
$ws.Cells.Item(6, 10).Value = 0.4415672215387696
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 2.642196
$ws.Cells.Item(6, 14).Value = 7.926588000000001
$ws.Cells.Item(6, 15).Value = 0.26568831615543
$ws.Cells.Item(6, 16).Value = 0.26568831615543
$ws.Cells.Item(6, 17).Value = 2.571365771096001
$ws.Cells.Item(6, 18).Value = 23.142291939864
$ws.Cells.Item(6, 19).Value = 0.1173192515600674
$ws.Cells.Item(6, 20).Value = 0.1173192515600674

# Row 7: FAPs -> FAPs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "F2"
$ws.Cells.Item(7, 3).Value = "Gp1ba"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 0.9731926666666667
$ws.Cells.Item(7, 8).Value = 2.919578
$ws.Cells.Item(7, 9).Value = 0.4415672215387696
$ws.Cells.Item(7, 10).Value = 0.4415672215387696
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 4.666004666666667
$ws.Cells.Item(7, 14).Value = 13.998014
$ws.Cells.Item(7, 15).Value = 0.469194156323015
$ws.Cells.Item(7, 16).Value = 0.4691941563230151
$ws.Cells.Item(7, 17).Value = 4.540921524232445
$ws.Cells.Item(7, 18).Value = 40.868293718092
$ws.Cells.Item(7, 19).Value = 0.2071807599697809
$ws.Cells.Item(7, 20).Value = 0.2071807599697809

# Row 8: FAPs -> M2 (ligand F2, receptor Gp1ba)
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "F2"
$ws.Cells.Item(8, 3).Value = "Gp1ba"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 0.9731926666666667
$ws.Cells.Item(8, 8).Value = 2.919578
$ws.Cells.Item(8, 9).Value = 0.4415672215387696
$ws.Cells.Item(8, 10).Value = 0.4415672215387696
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 1.114591666666667
$ws.Cells.Item(8, 14).Value = 3.343775
$ws.Cells.Item(8, 15).Value = 0.1120787341732184
$ws.Cells.Item(8, 16).Value = 0.1120787341732184
$ws.Cells.Item(8, 17).Value = 1.084712436327778
$ws.Cells.Item(8, 18).Value = 9.76241192695
$ws.Cells.Item(8, 19).Value = 0.0494902952424504
$ws.Cells.Item(8, 20).Value = 0.0494902952424504

# Row 9: FAPs -> sCs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "F2"
$ws.Cells.Item(9, 3).Value = "Gp1ba"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 0.9731926666666667
$ws.Cells.Item(9, 8).Value = 2.919578
$ws.Cells.Item(9, 9).Value = 0.4415672215387696
$ws.Cells.Item(9, 10).Value = 0.4415672215387696
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 1.521928
$ws.Cells.Item(9, 14).Value = 4.565784000000001
$ws.Cells.Item(9, 15).Value = 0.1530387933483365
$ws.Cells.Item(9, 16).Value = 0.1530387933483365
$ws.Cells.Item(9, 17).Value = 1.481129168794667
$ws.Cells.Item(9, 18).Value = 13.330162519152
$ws.Cells.Item(9, 19).Value = 0.06757691476647089
$ws.Cells.Item(9, 20).Value = 0.06757691476647089

# Row 10: M2 -> ECs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "F2"
$ws.Cells.Item(10, 3).Value = "Gp1ba"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1.0
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.048595
$ws.Cells.Item(10, 8).Value = 0.145785
$ws.Cells.Item(10, 9).Value = 0.02204903496054208
$ws.Cells.Item(10, 10).Value = 0.02204903496054208
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 2.642196
$ws.Cells.Item(10, 14).Value = 7.926588000000001
$ws.Cells.Item(10, 15).Value = 0.26568831615543
$ws.Cells.Item(10, 16).Value = 0.26568831615543
$ws.Cells.Item(10, 17).Value = 0.12839751462
$ws.Cells.Item(10, 18).Value = 1.15557763158
$ws.Cells.Item(10, 19).Value = 0.005858170971518632
$ws.Cells.Item(10, 20).Value = 0.005858170971518634

# Row 11: M2 -> FAPs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "F2"
$ws.Cells.Item(11, 3).Value = "Gp1ba"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 1.0
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.048595
$ws.Cells.Item(11, 8).Value = 0.145785
$ws.Cells.Item(11, 9).Value = 0.02204903496054208
$ws.Cells.Item(11, 10).Value = 0.02204903496054208
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 4.666004666666667
$ws.Cells.Item(11, 14).Value = 13.998014
$ws.Cells.Item(11, 15).Value = 0.469194156323015
$ws.Cells.Item(11, 16).Value = 0.4691941563230151
$ws.Cells.Item(11, 17).Value = 0.2267444967766667
$ws.Cells.Item(11, 18).Value = 2.04070047099
$ws.Cells.Item(11, 19).Value = 0.0103452783560482
$ws.Cells.Item(11, 20).Value = 0.0103452783560482

# Row 12: M2 -> M2 (ligand F2, receptor Gp1ba)
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "F2"
$ws.Cells.Item(12, 3).Value = "Gp1ba"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 1.0
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.048595
$ws.Cells.Item(12, 8).Value = 0.145785
$ws.Cells.Item(12, 9).Value = 0.02204903496054208
$ws.Cells.Item(12, 10).Value = 0.02204903496054208
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 1.114591666666667
$ws.Cells.Item(12, 14).Value = 3.343775
$ws.Cells.Item(12, 15).Value = 0.1120787341732184
$ws.Cells.Item(12, 16).Value = 0.1120787341732184
$ws.Cells.Item(12, 17).Value = 0.05416358204166666
$ws.Cells.Item(12, 18).Value = 0.487472238375
$ws.Cells.Item(12, 19).Value = 0.002471227928118594
$ws.Cells.Item(12, 20).Value = 0.002471227928118594

# Row 13: M2 -> sCs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "F2"
$ws.Cells.Item(13, 3).Value = "Gp1ba"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 1.0
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.048595
$ws.Cells.Item(13, 8).Value = 0.145785
$ws.Cells.Item(13, 9).Value = 0.02204903496054208
$ws.Cells.Item(13, 10).Value = 0.02204903496054208
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 1.521928
$ws.Cells.Item(13, 14).Value = 4.565784000000001
$ws.Cells.Item(13, 15).Value = 0.1530387933483365
$ws.Cells.Item(13, 16).Value = 0.1530387933483365
$ws.Cells.Item(13, 17).Value = 0.07395809116
$ws.Cells.Item(13, 18).Value = 0.6656228204400001
$ws.Cells.Item(13, 19).Value = 0.003374357704856646
$ws.Cells.Item(13, 20).Value = 0.003374357704856646

# Row 14: sCs -> ECs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "F2"
$ws.Cells.Item(14, 3).Value = "Gp1ba"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3.0
$ws.Cells.Item(14, 6).Value = 1.0
$ws.Cells.Item(14, 7).Value = 0.508125
$ws.Cells.Item(14, 8).Value = 1.524375
$ws.Cells.Item(14, 9).Value = 0.2305518240420917
$ws.Cells.Item(14, 10).Value = 0.2305518240420917
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 2.642196
$ws.Cells.Item(14, 14).Value = 7.926588000000001
$ws.Cells.Item(14, 15).Value = 0.26568831615543
$ws.Cells.Item(14, 16).Value = 0.26568831615543
$ws.Cells.Item(14, 17).Value = 1.3425658425
$ws.Cells.Item(14, 18).Value = 12.0830925825
$ws.Cells.Item(14, 19).Value = 0.06125492591630632
$ws.Cells.Item(14, 20).Value = 0.06125492591630633

# Row 15: sCs -> FAPs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "F2"
$ws.Cells.Item(15, 3).Value = "Gp1ba"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3.0
$ws.Cells.Item(15, 6).Value = 1.0
$ws.Cells.Item(15, 7).Value = 0.508125
$ws.Cells.Item(15, 8).Value = 1.524375
$ws.Cells.Item(15, 9).Value = 0.2305518240420917
$ws.Cells.Item(15, 10).Value = 0.2305518240420917
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 12).Value = 1.0
$ws.Cells.Item(15, 13).Value = 4.666004666666667
$ws.Cells.Item(15, 14).Value = 13.998014
$ws.Cells.Item(15, 15).Value = 0.469194156323015
$ws.Cells.Item(15, 16).Value = 0.4691941563230151
$ws.Cells.Item(15, 17).Value = 2.37091362125
$ws.Cells.Item(15, 18).Value = 21.33822259125
$ws.Cells.Item(15, 19).Value = 0.1081735685701614
$ws.Cells.Item(15, 20).Value = 0.1081735685701614

# Row 16: sCs -> M2 (ligand F2, receptor Gp1ba)
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "F2"
$ws.Cells.Item(16, 3).Value = "Gp1ba"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3.0
$ws.Cells.Item(16, 6).Value = 1.0
$ws.Cells.Item(16, 7).Value = 0.508125
$ws.Cells.Item(16, 8).Value = 1.524375
$ws.Cells.Item(16, 9).Value = 0.2305518240420917
$ws.Cells.Item(16, 10).Value = 0.2305518240420917
$ws.Cells.Item(16, 11).Value = 3.0
$ws.Cells.Item(16, 12).Value = 1.0
$ws.Cells.Item(16, 13).Value = 1.114591666666667
$ws.Cells.Item(16, 14).Value = 3.343775
$ws.Cells.Item(16, 15).Value = 0.1120787341732184
$ws.Cells.Item(16, 16).Value = 0.1120787341732184
$ws.Cells.Item(16, 17).Value = 0.566351890625
$ws.Cells.Item(16, 18).Value = 5.097167015625
$ws.Cells.Item(16, 19).Value = 0.02583995659996421
$ws.Cells.Item(16, 20).Value = 0.02583995659996421

# Row 17: sCs -> sCs (ligand F2, receptor Gp1ba)
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "F2"
$ws.Cells.Item(17, 3).Value = "Gp1ba"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3.0
$ws.Cells.Item(17, 6).Value = 1.0
$ws.Cells.Item(17, 7).Value = 0.508125
$ws.Cells.Item(17, 8).Value = 1.524375
$ws.Cells.Item(17, 9).Value = 0.2305518240420917
$ws.Cells.Item(17, 10).Value = 0.2305518240420917
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 12).Value = 1.0
$ws.Cells.Item(17, 13).Value = 1.521928
$ws.Cells.Item(17, 14).Value = 4.565784000000001
$ws.Cells.Item(17, 15).Value = 0.1530387933483365
$ws.Cells.Item(17, 16).Value = 0.1530387933483365
$ws.Cells.Item(17, 17).Value = 0.7733296650000001
$ws.Cells.Item(17, 18).Value = 6.959966985000001
$ws.Cells.Item(17, 19).Value = 0.03528337295565971
$ws.Cells.Item(17, 20).Value = 0.03528337295565971
